$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 468 (shifts existing rows 468:574 down to 469:575)
$ws.Rows(468).Insert()

# Populate the newly inserted row 468 with its data
$ws.Cells.Item(468, 1).Value = 5
$ws.Cells.Item(468, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(468, 3).Value = "Maule"
$ws.Cells.Item(468, 4).Value = 45244
$ws.Cells.Item(468, 5).Value = 7
$ws.Cells.Item(468, 6).Value = 100112006
$ws.Cells.Item(468, 7).Value = "Repollo"
$ws.Cells.Item(468, 8).Value = "Crespo record"
$ws.Cells.Item(468, 9).Value = "Primera"
$ws.Cells.Item(468, 10).Value = 3000
$ws.Cells.Item(468, 11).Value = 1200
$ws.Cells.Item(468, 12).Value = 1200
$ws.Cells.Item(468, 13).Value = 1200
$ws.Cells.Item(468, 14).Value = "$/unidad"
$ws.Cells.Item(468, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(468, 16).Value = 1200
$ws.Cells.Item(468, 17).Value = 1
$ws.Cells.Item(468, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date-number format used by the rest of column D
$ws.Cells.Item(468, 4).NumberFormat = $ws.Cells.Item(469, 4).NumberFormat
